# Apply the "DH with heatpump and boiler functional" changes to
# db_defaultEnergyAssets.xlsx
#
# Summary of data changes:
#  - conversionAssets sheet (3rd tab): row 10 becomes a new "DH_boiler_L" /
#    BOILER asset; the former row 10 (methane furnace) moves down to row 11,
#    and the former row 11 (hydrogen furnace) moves down to a brand new
#    row 12.
#  - storageAssets sheet (4th tab): row 11 (District_Heating_heat_buffer_HT_S)
#    ambientTempType changes from AIR to GROUND.
#  - The active/selected sheet & cell selection move from conversionAssets
#    to storageAssets.

$wb = $excel.ActiveWorkbook

$wsConversion = $wb.Worksheets.Item("conversionAssets")
$wsStorage    = $wb.Worksheets.Item("storageAssets")

# --- conversionAssets: insert the new DH_boiler_L / BOILER row at row 10,
#     pushing the furnace rows down by one ---

# Former row 11 (Industrial_hydrogen_furnace / HYDROGEN_FURNACE) -> row 12
$wsConversion.Range("A12").Value = 11
$wsConversion.Range("B12").Value = "Industrial_hydrogen_furnace"
$wsConversion.Range("C12").Value = "CONVERSION"
$wsConversion.Range("D12").Value = "HYDROGEN_FURNACE"
$wsConversion.Range("E12").Value = 0
$wsConversion.Range("F12").Value = 300
$wsConversion.Range("G12").Value = 0.99
$wsConversion.Range("H12").Value = 120

# Former row 10 (Industrial_methane_furnace / METHANE_FURNACE) -> row 11
$wsConversion.Range("B11").Value = "Industrial_methane_furnace"
$wsConversion.Range("D11").Value = "METHANE_FURNACE"

# New row 10: DH_boiler_L / BOILER
$wsConversion.Range("B10").Value = "DH_boiler_L"
$wsConversion.Range("D10").Value = "BOILER"
$wsConversion.Range("E10").Value = 300
$wsConversion.Range("F10").Value = 0
$wsConversion.Range("H10").Value = 100

# --- storageAssets: District_Heating_heat_buffer_HT_S ambientTempType
#     AIR -> GROUND ---
$wsStorage.Range("N11").Value = "GROUND"

# --- Selection on conversionAssets ends up parked on H10 (last touched
#     cell), while the workbook's active sheet/selection moves to
#     storageAssets!G11 ---
$wsConversion.Range("H10").Select()
$wsStorage.Activate()
$wsStorage.Range("G11").Select()
